$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New block: "use sentiment word frequency + word counts:" (rows 14-16)
# Written first so it mints a brand new shared-string entry rather than
# reusing/mutating the "learning rate decay" one below.
$ws.Range("A14").Value = "use sentiment word frequency + word counts:"
$ws.Range("B14").Value = "epoch1"
$ws.Range("C14").Value = 0.61737329228550897
$ws.Range("B15").Value = "epoch2"
$ws.Range("C15").Value = 0.62109860854033605
$ws.Range("B16").Value = "epoch3"
$ws.Range("C16").Value = 0.62598893967624802

# Fix the typo in the existing "learning rate decay" label (row 10, column A)
$ws.Range("A10").Value = "use sentiment word frequency + learning rate decay:"

# Apply the same bold-ish header style used on A2/A6/A10 to the new A14 label
$ws.Range("A10").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New D6 value (mirrors D2/D10 style extra metric column)
$ws.Range("D6").Value = 0.70199999999999996

# New block: extra raw values (rows 18-20), column C only
$ws.Range("C18").Value = 0.60709525571029099
$ws.Range("C19").Value = 0.627169190381766
$ws.Range("C20").Value = 0.61413615973564695

# Update the sheet view to match the new scroll/selection state
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("E17").Select()
